# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> linked from the notes master (was "Office Theme")
#   ppt/theme/theme2.xml -> linked from the slide master  (was "Integral" / "Red Violet")
#
# The edit swaps the two themes' contents: the slide master's theme becomes
# the plain "Office" palette, and the notes-master theme becomes the
# "Red Violet" / Integral palette. Font scheme and format scheme are already
# identical between the two themes ("Office"), so only the twelve theme
# colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) need to
# change on the slide master's theme.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Target palette = the former "Office Theme" colour scheme (was theme1.xml,
# now becomes the slide master's live theme).
$colors.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$colors.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$colors.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
